# Review Tracker update: add Amr's reviews (rows 19-22) to the "Review" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: Review ---
$ws.Range("B19").Value = "Review "
$ws.Range("C19").Value = "Each task from the team , coach and Eng/ moahmed hassan should has review."
$ws.Range("D19").Value = "Amr"
$ws.Range("F19").Value = 45082

# --- Row 20: RTM ---
$ws.Range("B20").Value = "RTM"
$ws.Range("C20").Value = "1. update all use cases.`n2. Add alternative  flows from srs."
$ws.Range("D20").Value = "Amr"
$ws.Range("F20").Value = 45082

# --- Row 21: Design Document ---
$ws.Range("B21").Value = "Design Document"
$ws.Range("C21").Value = "Add ids to be reahable."
$ws.Range("D21").Value = "Amr"
$ws.Range("F21").Value = 45082

# --- Row 22: ERD ---
$ws.Range("B22").Value = "ERD"
$ws.Range("C22").Value = "change ""manage"" relation to be add/delete."
$ws.Range("D22").Value = "Amr"
$ws.Range("F22").Value = 45082

# Match the style used by C18 (center/middle, wrap) for the new comment cells,
# and the plain content style (style used by B/D/E column cells elsewhere, e.g. B14/D14) for B/D/E.
$commentStyle = $ws.Range("C18").Style
$ws.Range("C19:C22").Style = $commentStyle

$plainStyle = $ws.Range("B14").Style
$ws.Range("B19:B22").Style = $plainStyle
$ws.Range("D19:D22").Style = $plainStyle
$ws.Range("E19:E22").Style = $plainStyle

$dateStyle = $ws.Range("F14").Style
$ws.Range("F19:F22").Style = $dateStyle

# Clear any stray value carried from style copy (E column is empty in the target rows)
$ws.Range("E19:E22").ClearContents()

# Update the frozen-pane anchor and the active selection to match the latest edit location.
$ws.Application.ActiveWindow.Panes.Item(2).ScrollRow = 17
$ws.Range("G19").Select()
